# Hortaliza, Feria Lagunitas de Puerto Montt - Pepino ensalada
# A new daily price-report row is inserted at row 218 (pushing the
# existing rows 218-279 down to 219-280); the rest of the sheet is
# unchanged aside from the shift and the growth of the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 218; everything at/after 218 shifts down
# by one (old row 218 -> new row 219, ..., old row 279 -> new row 280).
$ws.Rows.Item(218).Insert()

# Populate the newly inserted row 218 with the new observation.
$ws.Range("A218").Value = 4
$ws.Range("B218").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C218").Value = "Los Lagos"
$ws.Range("D218").Value = 44736
$ws.Range("E218").Value = 10
$ws.Range("F218").Value = 100112043
$ws.Range("G218").Value = "Pepino ensalada"
$ws.Range("H218").Value = "Sin especificar"
$ws.Range("I218").Value = "Primera"
$ws.Range("J218").Value = 400
$ws.Range("K218").Value = 23000
$ws.Range("L218").Value = 23000
$ws.Range("M218").Value = 23000
$ws.Range("N218").Value = "`$/caja 60 unidades"
$ws.Range("O218").Value = "Región de Arica y Parinacota"
$ws.Range("P218").Value = 383
$ws.Range("Q218").Value = 60
$ws.Range("R218").Value = "Hortaliza"
